$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Prófa"
$ws.Range("A3").Value = "Rejpal"
$ws.Range("A4").Value = "Štístko"
$ws.Range("A5").Value = "Dřímal"
$ws.Range("A6").Value = "Stydlín"
$ws.Range("A7").Value = "Kejchal"
$ws.Range("A8").Value = "Šmudla"

$ws.Range("A1").Select()
